# Fruta / hortaliza, semanal
# Insert two new weekly price rows (for 2023-03-30) above the existing
# "Terminal Hortofrutícola Agro Chillán - Ciruela / Angeleno" records,
# pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 88-89; Excel copies formatting (e.g. the date
# style in column D) from the row above, and the sheet dimension/row
# references below automatically shift down by two.
$ws.Rows("88:89").Insert()

# --- New row 88: Angeleno / Especial --------------------------------
$ws.Cells.Item(88, 1).Value = 7
$ws.Cells.Item(88, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(88, 3).Value = "Ñuble"
$ws.Cells.Item(88, 4).Value = 45015
$ws.Cells.Item(88, 5).Value = 16
$ws.Cells.Item(88, 6).Value = "Fruta"
$ws.Cells.Item(88, 7).Value = 100103
$ws.Cells.Item(88, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(88, 9).Value = 100103002
$ws.Cells.Item(88, 10).Value = "Ciruela"
$ws.Cells.Item(88, 11).Value = "Angeleno"
$ws.Cells.Item(88, 12).Value = "Especial"
$ws.Cells.Item(88, 13).Value = 50
$ws.Cells.Item(88, 14).Value = 12000
$ws.Cells.Item(88, 15).Value = 12000
$ws.Cells.Item(88, 16).Value = 12000
$ws.Cells.Item(88, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(88, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(88, 19).Value = 667
$ws.Cells.Item(88, 20).Value = 18

# --- New row 89: Angeleno / Primera ----------------------------------
$ws.Cells.Item(89, 1).Value = 7
$ws.Cells.Item(89, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(89, 3).Value = "Ñuble"
$ws.Cells.Item(89, 4).Value = 45015
$ws.Cells.Item(89, 5).Value = 16
$ws.Cells.Item(89, 6).Value = "Fruta"
$ws.Cells.Item(89, 7).Value = 100103
$ws.Cells.Item(89, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(89, 9).Value = 100103002
$ws.Cells.Item(89, 10).Value = "Ciruela"
$ws.Cells.Item(89, 11).Value = "Angeleno"
$ws.Cells.Item(89, 12).Value = "Primera"
$ws.Cells.Item(89, 13).Value = 40
$ws.Cells.Item(89, 14).Value = 11000
$ws.Cells.Item(89, 15).Value = 11000
$ws.Cells.Item(89, 16).Value = 11000
$ws.Cells.Item(89, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(89, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(89, 19).Value = 611
$ws.Cells.Item(89, 20).Value = 18

Write-Host "Inserted rows 88-89 and populated new weekly records."
